# Update market data / stock prices: insert a new (earlier) trade entry into the
# "Trading History" sheet, ahead of the existing row, shifting the existing
# trade row down from row 5 to row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trading History")

$dateFormat = "yyyy-mm-dd h:mm:ss"

# Move the existing trade (previously in row 5) down to row 6.
$ws.Range("A6").Value = 46063
$ws.Range("A6").NumberFormat = $dateFormat
$ws.Range("B6").Value = "NSE"
$ws.Range("C6").Value = "Buy"
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 890.01
$ws.Range("F6").Value = 4472.3
$ws.Range("G6").Value = "CN#252611730667"
$ws.Range("I6").Value = 22.25
$ws.Range("J6").Formula = "=Index!`$C`$2"

# Add the new, earlier trade entry into row 5.
$ws.Range("A5").Value = 46062
$ws.Range("A5").NumberFormat = $dateFormat
$ws.Range("B5").Value = "NSE"
$ws.Range("C5").Value = "Buy"
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 891.95
$ws.Range("F5").Value = 4482.05
$ws.Range("G5").Value = "CN#252611665409"
$ws.Range("I5").Value = 22.3
$ws.Range("J5").Formula = "=Index!`$C`$2"
